{"js": "// Append \". \" to the \"Proof of vehicle registration...\" paragraph, then\n// insert a new set of numbered-requirement paragraphs (with blank-line\n// spacers) right after the existing empty paragraph that follows it.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"Proof of vehicle registration...\" paragraph by its text.\nlet proofPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Proof of vehicle registration is essential\") !== -1) {\n    proofPara = paragraphs.items[i];\n    break;\n  }\n}\nif (!proofPara) {\n  throw new Error(\"Could not find 'Proof of vehicle registration...' paragraph\");\n}\n\n// Append a trailing \". \" to that paragraph's text.\nproofPara.insertText(\". \", \"End\");\nawait context.sync();\n\n// The blank paragraph right after it is where the new content begins.\nlet anchor = proofPara.getNext();\nanchor.load(\"text\");\nawait context.sync();\n\n// New paragraphs to insert, in order. `null` denotes a blank paragraph.\nconst newParagraphs = [\n  \"1. Valid Registration\",\n  \"The vehicle must have valid, current registration in the state where the incident occurred at the time of the incident in order to successfully process the claim.\",\n  null,\n  \"Details:\",\n  \"- Registration must be current and not expired at the time of incident\",\n  \"- Registration must be to the same state as where incident occurred \",\n  \"- If registration is expired or not valid for the incident state at the time of occurrence, the claim may be denied\",\n  null,\n  \"2. Registration in Policyholder's Name \",\n  \"The vehicle must be registered in the name of the policyholder filing the insurance claim. \",\n  null,\n  \"Details: \",\n  \"- The name on the vehicle registration must match the name of the policyholder\",\n  \"- If the registered owner differs from the policyholder, additional documentation connecting the policyholder to the registered owner may be required to process the claim\",\n  null,\n  \"3. Provide Registration Documentation\",\n  \"The policyholder must provide documentation supporting the vehicle's registration status along with the claim submission.\",\n  null,\n  \"Details:\",\n  \"- Acceptable docs: vehicle registration card, DMV registration receipt showing effective dates and name, vehicle title showing registration in effect at time of incident\",\n  \"- Documentation must show effective registration status as described in Requirement 1\",\n  \"- Failure to provide documentation may delay processing or result in claim denial\",\n  null,\n  \"4. Out of State Registration Situations \",\n  \"If the vehicle has out of state registration plates compared to where the incident occurred, the out of state registration status must still meet all other requirements. \",\n  null,\n  \"Details: \",\n  \"- Out of state registration must still be valid  and in policyholder's name\",\n  \"- Applicable state's registration rules apply (where incident happened)\",\n  \"- May require additional investigation into registration status across multiple states\",\n];\n\nfor (const text of newParagraphs) {\n  anchor = anchor.insertParagraph(text === null ? \"\" : text, \"After\");\n  await context.sync();\n}\n", "ps1": "# Append \". \" to the \"Proof of vehicle registration...\" paragraph, then\n# insert a new set of numbered-requirement paragraphs (with blank-line\n# spacers) right after the existing empty paragraph that follows it.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Proof of vehicle registration...\" paragraph by its text.\n$proofPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"Proof of vehicle registration is essential*\") {\n        $proofPara = $p\n        break\n    }\n}\nif ($proofPara -eq $null) {\n    throw \"Could not find 'Proof of vehicle registration...' paragraph\"\n}\n\n# Append a trailing \". \" to that paragraph's text.\n$proofRange = $proofPara.Range\n$proofRange.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark\n$proofRange.InsertAfter(\". \")\n\n# The blank paragraph right after it is where the new content begins.\n$anchor = $proofPara.Next().Range\n$anchor.Collapse(0)\n\n# New paragraphs to insert, in order. \"\" denotes a blank paragraph.\n$newParagraphs = @(\n    \"1. Valid Registration\",\n    \"The vehicle must have valid, current registration in the state where the incident occurred at the time of the incident in order to successfully process the claim.\",\n    \"\",\n    \"Details:\",\n    \"- Registration must be current and not expired at the time of incident\",\n    \"- Registration must be to the same state as where incident occurred \",\n    \"- If registration is expired or not valid for the incident state at the time of occurrence, the claim may be denied\",\n    \"\",\n    \"2. Registration in Policyholder's Name \",\n    \"The vehicle must be registered in the name of the policyholder filing the insurance claim. \",\n    \"\",\n    \"Details: \",\n    \"- The name on the vehicle registration must match the name of the policyholder\",\n    \"- If the registered owner differs from the policyholder, additional documentation connecting the policyholder to the registered owner may be required to process the claim\",\n    \"\",\n    \"3. Provide Registration Documentation\",\n    \"The policyholder must provide documentation supporting the vehicle's registration status along with the claim submission.\",\n    \"\",\n    \"Details:\",\n    \"- Acceptable docs: vehicle registration card, DMV registration receipt showing effective dates and name, vehicle title showing registration in effect at time of incident\",\n    \"- Documentation must show effective registration status as described in Requirement 1\",\n    \"- Failure to provide documentation may delay processing or result in claim denial\",\n    \"\",\n    \"4. Out of State Registration Situations \",\n    \"If the vehicle has out of state registration plates compared to where the incident occurred, the out of state registration status must still meet all other requirements. \",\n    \"\",\n    \"Details: \",\n    \"- Out of state registration must still be valid  and in policyholder's name\",\n    \"- Applicable state's registration rules apply (where incident happened)\",\n    \"- May require additional investigation into registration status across multiple states\"\n)\n\nforeach ($t in $newParagraphs) {\n    $anchor.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Last\n    if ($t -ne \"\") {\n        $newPara.Range.InsertBefore($t)\n    }\n    $anchor = $d.Paragraphs.Last.Range\n    $anchor.Collapse(0)\n}\n"}
